$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.824.74"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "1.640.29"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.17%  "

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "216.19"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("E10").Value = "  -1.34%  "

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "4.27"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").Value = "1.866.78"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").Value = "1.642.33"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("E16").Value = "  +0.29%  "

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "63.07"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "25.863.60"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("E20").Value = "  +2.32%  "

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "192.91"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.51%  "

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "9.97"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +0.62%  "

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "6.35"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +2.83%  "

$ws.Range("E24").Value = "  +5.73%  "

$ws.Range("E25").Value = "  -0.10%  "

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "142.39"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +2.66%  "

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.124"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +0.97%  "

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "6.94"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +1.93%  "

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "15.54"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").Value = "  +0.59%  "

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "2.38"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -0.33%  "

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.908"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").Value = "1.135.30"
$ws.Range("E37").Value = "  +1.27%  "

$ws.Range("E38").Value = "  -1.59%  "

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.547"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("E41").Value = "  +0.09%  "

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "5.60"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +1.72%  "

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "100.74"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +1.30%  "

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.807"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("D45").Value = "1.775.80"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("E46").Value = "  -0.22%  "

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "55.39"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  +6.73%  "

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.418"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("E51").Value = "  +2.58%  "
